{"js": "// Update the worksheet date and each \"NNN\u00f7N=\" division prompt to the\n// new values from the regenerated output (commit c986bee).\nconst replacements = [\n  [\"2024-10-18 Friday\", \"2024-10-19 Saturday\"],\n  [\"885\u00f74=\", \"151\u00f79=\"],\n  [\"594\u00f78=\", \"969\u00f79=\"],\n  [\"447\u00f77=\", \"301\u00f72=\"],\n  [\"360\u00f76=\", \"497\u00f78=\"],\n  [\"979\u00f75=\", \"186\u00f79=\"],\n  [\"174\u00f75=\", \"393\u00f78=\"],\n  [\"303\u00f74=\", \"618\u00f74=\"],\n  [\"550\u00f79=\", \"356\u00f75=\"],\n  [\"704\u00f79=\", \"478\u00f76=\"],\n  [\"550\u00f77=\", \"882\u00f78=\"],\n  [\"138\u00f73=\", \"651\u00f72=\"],\n  [\"387\u00f74=\", \"458\u00f75=\"],\n  [\"346\u00f76=\", \"810\u00f77=\"],\n  [\"440\u00f75=\", \"817\u00f75=\"],\n  [\"138\u00f78=\", \"402\u00f79=\"],\n  [\"920\u00f78=\", \"166\u00f79=\"],\n  [\"657\u00f74=\", \"804\u00f78=\"],\n  [\"149\u00f76=\", \"926\u00f77=\"],\n  [\"382\u00f76=\", \"588\u00f74=\"],\n  [\"994\u00f79=\", \"756\u00f72=\"],\n  [\"473\u00f74=\", \"215\u00f73=\"],\n  [\"382\u00f73=\", \"939\u00f79=\"],\n  [\"653\u00f76=\", \"467\u00f78=\"],\n  [\"624\u00f73=\", \"887\u00f78=\"],\n  [\"840\u00f74=\", \"692\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and each \"NNN\u00f7N=\" division prompt to the\n# new values from the regenerated output (commit c986bee).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-18 Friday\", \"2024-10-19 Saturday\"),\n    @(\"885\u00f74=\", \"151\u00f79=\"),\n    @(\"594\u00f78=\", \"969\u00f79=\"),\n    @(\"447\u00f77=\", \"301\u00f72=\"),\n    @(\"360\u00f76=\", \"497\u00f78=\"),\n    @(\"979\u00f75=\", \"186\u00f79=\"),\n    @(\"174\u00f75=\", \"393\u00f78=\"),\n    @(\"303\u00f74=\", \"618\u00f74=\"),\n    @(\"550\u00f79=\", \"356\u00f75=\"),\n    @(\"704\u00f79=\", \"478\u00f76=\"),\n    @(\"550\u00f77=\", \"882\u00f78=\"),\n    @(\"138\u00f73=\", \"651\u00f72=\"),\n    @(\"387\u00f74=\", \"458\u00f75=\"),\n    @(\"346\u00f76=\", \"810\u00f77=\"),\n    @(\"440\u00f75=\", \"817\u00f75=\"),\n    @(\"138\u00f78=\", \"402\u00f79=\"),\n    @(\"920\u00f78=\", \"166\u00f79=\"),\n    @(\"657\u00f74=\", \"804\u00f78=\"),\n    @(\"149\u00f76=\", \"926\u00f77=\"),\n    @(\"382\u00f76=\", \"588\u00f74=\"),\n    @(\"994\u00f79=\", \"756\u00f72=\"),\n    @(\"473\u00f74=\", \"215\u00f73=\"),\n    @(\"382\u00f73=\", \"939\u00f79=\"),\n    @(\"653\u00f76=\", \"467\u00f78=\"),\n    @(\"624\u00f73=\", \"887\u00f78=\"),\n    @(\"840\u00f74=\", \"692\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
